{"js": "// Prepare public release and harden processing reliability.\n//\n// The document starts with a title, an intro paragraph, a \"Key/Value\"\n// table, an \"Here is another table:\" paragraph and a Product/Price/Stock\n// table with a fruit row set (incl. an Orange row). The target state:\n//   - Title text -> \"Test Document with Table\"\n//   - Intro paragraph text -> \"This is a test document.\"\n//   - The first (Key/Value) table and the \"Here is another table:\"\n//     paragraph are removed entirely.\n//   - The remaining table gets the \"LightGrid-Accent1\" table style and its\n//     contents become a Name/Age/City people table (Alice/Bob), with the\n//     trailing Orange-equivalent (3rd data) row removed.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nbody.tables.load(\"items\");\nawait context.sync();\n\n// Title + intro paragraph (always the first two body-level paragraphs,\n// ahead of anything living inside a table cell).\nconst paras = body.paragraphs.items;\nparas[0].insertText(\"Test Document with Table\", \"Replace\");\nparas[1].insertText(\"This is a test document.\", \"Replace\");\n\nconst tables = body.tables.items;\nconst keyValueTable = tables[0];\nconst productTable = tables[1];\n\n// Drop the \"Here is another table:\" lead-in paragraph and the whole\n// Key/Value table.\nconst leadInParagraph = productTable.getParagraphBefore();\nleadInParagraph.delete();\nkeyValueTable.delete();\nawait context.sync();\n\n// Re-fetch: the surviving table's anchor shifted once the Key/Value table\n// (and the paragraph ahead of it) were removed from the body.\nbody.tables.load(\"items\");\nawait context.sync();\nconst table = body.tables.items[0];\n\ntable.style = \"LightGrid-Accent1\";\n\ntable.rows.load(\"items\");\nawait context.sync();\nconst rows = table.rows.items;\n\n// Header row: Product/Price/Stock -> Name/Age/City\nrows[0].getCell(0).body.insertText(\"Name\", \"Replace\");\nrows[0].getCell(1).body.insertText(\"Age\", \"Replace\");\nrows[0].getCell(2).body.insertText(\"City\", \"Replace\");\n\n// Row 1: Apple/1.99/100 -> Alice/30/NYC\nrows[1].getCell(0).body.insertText(\"Alice\", \"Replace\");\nrows[1].getCell(1).body.insertText(\"30\", \"Replace\");\nrows[1].getCell(2).body.insertText(\"NYC\", \"Replace\");\n\n// Row 2: Banana/0.99/50 -> Bob/25/LA\nrows[2].getCell(0).body.insertText(\"Bob\", \"Replace\");\nrows[2].getCell(1).body.insertText(\"25\", \"Replace\");\nrows[2].getCell(2).body.insertText(\"LA\", \"Replace\");\n\n// Row 3 (Orange/2.49/75) is dropped entirely.\nrows[3].delete();\n\nawait context.sync();\n", "ps1": "# Prepare public release and harden processing reliability.\n#\n# The document starts with a title, an intro paragraph, a \"Key/Value\"\n# table, an \"Here is another table:\" paragraph and a Product/Price/Stock\n# table with a fruit row set (incl. an Orange row). Target state:\n#   - Title text -> \"Test Document with Table\"\n#   - Intro paragraph text -> \"This is a test document.\"\n#   - The first (Key/Value) table and the \"Here is another table:\"\n#     paragraph are removed entirely.\n#   - The remaining table gets the \"LightGrid-Accent1\" table style and its\n#     contents become a Name/Age/City people table (Alice/Bob), with the\n#     trailing Orange-equivalent (3rd data) row removed.\n\n$d = $word.ActiveDocument\n\n# Title + intro paragraph (always the first two paragraphs of the story).\n$d.Paragraphs.Item(1).Range.Text = \"Test Document with Table\"\n$d.Paragraphs.Item(2).Range.Text = \"This is a test document.\"\n\n# Locate the \"Here is another table:\" lead-in paragraph (Range.Text carries\n# a trailing paragraph mark, so compare trimmed) and remove it completely.\n$idx = 0\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.Trim() -eq \"Here is another table:\") {\n        $idx = $i\n    }\n}\n$d.Paragraphs.Item($idx).Range.Delete()\n\n# Remove the whole Key/Value table (the first table in the document).\n$d.Tables.Item(1).Delete()\n\n# The surviving table is now Tables.Item(1); re-style it and rewrite its\n# contents into the Name/Age/City roster.\n$t = $d.Tables.Item(1)\n$t.Style = \"LightGrid-Accent1\"\n\n$t.Cell(1, 1).Range.Text = \"Name\"\n$t.Cell(1, 2).Range.Text = \"Age\"\n$t.Cell(1, 3).Range.Text = \"City\"\n\n$t.Cell(2, 1).Range.Text = \"Alice\"\n$t.Cell(2, 2).Range.Text = \"30\"\n$t.Cell(2, 3).Range.Text = \"NYC\"\n\n$t.Cell(3, 1).Range.Text = \"Bob\"\n$t.Cell(3, 2).Range.Text = \"25\"\n$t.Cell(3, 3).Range.Text = \"LA\"\n\n# Drop the trailing Orange-equivalent (4th) row entirely.\n$t.Rows.Item(4).Delete()\n"}
